$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 639
$ws.Cells.Item(639,1).Value2 = 45192.67608010417
$ws.Cells.Item(639,2).Value = "mjh8429@naver.com"
$ws.Cells.Item(639,3).Value = "디지털미디어콘텐츠"
$ws.Cells.Item(639,4).Value2 = 20192535
$ws.Cells.Item(639,5).Value = "민지혜"
$ws.Cells.Item(639,6).Value = "78:22"
$ws.Cells.Item(639,7).Value2 = 0.15
$ws.Cells.Item(639,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(639,9).Value = "779만 명"
$ws.Cells.Item(639,10).Value2 = 0.151
$ws.Cells.Item(639,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(639,12).Value = "Red"
$ws.Cells.Item(639,13).Value = "모름/무응답"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A639:L639")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M639")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

# Row 640
$ws.Cells.Item(640,1).Value2 = 45192.67764202546
$ws.Cells.Item(640,2).Value = "041030top@naver.com"
$ws.Cells.Item(640,3).Value = "데이터사이언스"
$ws.Cells.Item(640,4).Value2 = 20233257
$ws.Cells.Item(640,5).Value = "최영국"
$ws.Cells.Item(640,6).Value = "74:26"
$ws.Cells.Item(640,7).Value2 = 0.2
$ws.Cells.Item(640,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(640,9).Value = "952만 명"
$ws.Cells.Item(640,10).Value2 = 0.059
$ws.Cells.Item(640,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(640,12).Value = "Black"
$ws.Cells.Item(640,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A640:L640")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N640")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 641
$ws.Cells.Item(641,1).Value2 = 45192.68418751158
$ws.Cells.Item(641,2).Value = "ekgus0916@naver.com"
$ws.Cells.Item(641,3).Value = "체육학과"
$ws.Cells.Item(641,4).Value2 = 20217125
$ws.Cells.Item(641,5).Value = "김다현"
$ws.Cells.Item(641,6).Value = "74:26"
$ws.Cells.Item(641,7).Value2 = 0.1
$ws.Cells.Item(641,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(641,9).Value = "779만 명"
$ws.Cells.Item(641,10).Value2 = 0.151
$ws.Cells.Item(641,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(641,12).Value = "Red"
$ws.Cells.Item(641,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A641:L641")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M641")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

# Row 642
$ws.Cells.Item(642,1).Value2 = 45192.71004814815
$ws.Cells.Item(642,2).Value = "jinyoung05503@gmail.com"
$ws.Cells.Item(642,3).Value = "생명과학과"
$ws.Cells.Item(642,4).Value2 = 20233534
$ws.Cells.Item(642,5).Value = "이진영"
$ws.Cells.Item(642,6).Value = "75:25"
$ws.Cells.Item(642,7).Value2 = 0.2
$ws.Cells.Item(642,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(642,9).Value = "779만 명"
$ws.Cells.Item(642,10).Value2 = 0.151
$ws.Cells.Item(642,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(642,12).Value = "Red"
$ws.Cells.Item(642,13).Value = "모름/무응답"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A642:L642")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M642")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

# Row 643
$ws.Cells.Item(643,1).Value2 = 45192.72864728009
$ws.Cells.Item(643,2).Value = "soccert71@naver.com"
$ws.Cells.Item(643,3).Value = "심리학과"
$ws.Cells.Item(643,4).Value2 = 20182125
$ws.Cells.Item(643,5).Value = "육정민"
$ws.Cells.Item(643,6).Value = "74:26"
$ws.Cells.Item(643,7).Value2 = 0.2
$ws.Cells.Item(643,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(643,9).Value = "952만 명"
$ws.Cells.Item(643,10).Value2 = 0.059
$ws.Cells.Item(643,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(643,12).Value = "Red"
$ws.Cells.Item(643,13).Value = "모름/무응답"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A643:L643")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M643")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

# Row 644
$ws.Cells.Item(644,1).Value2 = 45192.72900384259
$ws.Cells.Item(644,2).Value = "han7434370@naver.com"
$ws.Cells.Item(644,3).Value = "체육학과"
$ws.Cells.Item(644,4).Value2 = 20224152
$ws.Cells.Item(644,5).Value = "한진우"
$ws.Cells.Item(644,6).Value = "77:23"
$ws.Cells.Item(644,7).Value2 = 0.2
$ws.Cells.Item(644,8).Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Cells.Item(644,9).Value = "166만 명"
$ws.Cells.Item(644,10).Value2 = 0.151
$ws.Cells.Item(644,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(644,12).Value = "Red"
$ws.Cells.Item(644,13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A644:L644")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M644")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

# Row 645
$ws.Cells.Item(645,1).Value2 = 45192.73253275463
$ws.Cells.Item(645,2).Value = "milovany03@gmail.com"
$ws.Cells.Item(645,3).Value = "사회학과"
$ws.Cells.Item(645,4).Value2 = 20202223
$ws.Cells.Item(645,5).Value = "박진옥"
$ws.Cells.Item(645,6).Value = "74:26"
$ws.Cells.Item(645,7).Value2 = 0.2
$ws.Cells.Item(645,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(645,9).Value = "952만 명"
$ws.Cells.Item(645,10).Value2 = 0.059
$ws.Cells.Item(645,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(645,12).Value = "Black"
$ws.Cells.Item(645,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A645:L645")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N645")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 646
$ws.Cells.Item(646,1).Value2 = 45192.75874024305
$ws.Cells.Item(646,2).Value = "catboom5329@gmail.com"
$ws.Cells.Item(646,3).Value = "체육학과"
$ws.Cells.Item(646,4).Value2 = 20234121
$ws.Cells.Item(646,5).Value = "박주현"
$ws.Cells.Item(646,6).Value = "78:22"
$ws.Cells.Item(646,7).Value2 = 0.25
$ws.Cells.Item(646,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(646,9).Value = "952만 명"
$ws.Cells.Item(646,10).Value2 = 0.151
$ws.Cells.Item(646,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(646,12).Value = "Black"
$ws.Cells.Item(646,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A646:L646")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N646")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 647
$ws.Cells.Item(647,1).Value2 = 45192.7610846875
$ws.Cells.Item(647,2).Value = "yeshin05@naver.com"
$ws.Cells.Item(647,3).Value = "미래융합스쿨"
$ws.Cells.Item(647,4).Value2 = 20236639
$ws.Cells.Item(647,5).Value = "최예원"
$ws.Cells.Item(647,6).Value = "76:24"
$ws.Cells.Item(647,7).Value2 = 0.25
$ws.Cells.Item(647,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(647,9).Value = "779만 명"
$ws.Cells.Item(647,10).Value2 = 0.151
$ws.Cells.Item(647,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(647,12).Value = "Red"
$ws.Cells.Item(647,13).Value = "모름/무응답"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A647:L647")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M647")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

# Row 648
$ws.Cells.Item(648,1).Value2 = 45192.773353344906
$ws.Cells.Item(648,2).Value = "ind_b3@naver.com"
$ws.Cells.Item(648,3).Value = "미디어스쿨"
$ws.Cells.Item(648,4).Value2 = 20232523
$ws.Cells.Item(648,5).Value = "김지안"
$ws.Cells.Item(648,6).Value = "76:24"
$ws.Cells.Item(648,7).Value2 = 0.2
$ws.Cells.Item(648,8).Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Cells.Item(648,9).Value = "779만 명"
$ws.Cells.Item(648,10).Value2 = 0.059
$ws.Cells.Item(648,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(648,12).Value = "Black"
$ws.Cells.Item(648,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A648:L648")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N648")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 649
$ws.Cells.Item(649,1).Value2 = 45192.77901452546
$ws.Cells.Item(649,2).Value = "jamesjm0612@gmail.com"
$ws.Cells.Item(649,3).Value = "영어영문학과"
$ws.Cells.Item(649,4).Value2 = 20231231
$ws.Cells.Item(649,5).Value = "정재민"
$ws.Cells.Item(649,6).Value = "75:25"
$ws.Cells.Item(649,7).Value2 = 0.15
$ws.Cells.Item(649,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(649,9).Value = "166만 명"
$ws.Cells.Item(649,10).Value2 = 0.374
$ws.Cells.Item(649,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(649,12).Value = "Black"
$ws.Cells.Item(649,14).Value = "모름/무응답"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A649:L649")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N649")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 650
$ws.Cells.Item(650,1).Value2 = 45192.787771874995
$ws.Cells.Item(650,2).Value = "kkhe2370@naver.com"
$ws.Cells.Item(650,3).Value = "광고홍보학과"
$ws.Cells.Item(650,4).Value2 = 20202638
$ws.Cells.Item(650,5).Value = "전혜린"
$ws.Cells.Item(650,6).Value = "77:23"
$ws.Cells.Item(650,7).Value2 = 0.1
$ws.Cells.Item(650,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(650,9).Value = "166만 명"
$ws.Cells.Item(650,10).Value2 = 0.151
$ws.Cells.Item(650,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(650,12).Value = "Black"
$ws.Cells.Item(650,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A650:L650")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N650")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 651
$ws.Cells.Item(651,1).Value2 = 45192.790401087965
$ws.Cells.Item(651,2).Value = "jehuncho03@gmail.com"
$ws.Cells.Item(651,3).Value = "글로벌비즈니스"
$ws.Cells.Item(651,4).Value2 = 20226425
$ws.Cells.Item(651,5).Value = "조제헌"
$ws.Cells.Item(651,6).Value = "74:26"
$ws.Cells.Item(651,7).Value2 = 0.3
$ws.Cells.Item(651,8).Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Cells.Item(651,9).Value = "779만 명"
$ws.Cells.Item(651,10).Value2 = 0.059
$ws.Cells.Item(651,11).Value = "중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다"
$ws.Cells.Item(651,12).Value = "Black"
$ws.Cells.Item(651,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A651:L651")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N651")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 652
$ws.Cells.Item(652,1).Value2 = 45192.79356601852
$ws.Cells.Item(652,2).Value = "dms95123@naver.com"
$ws.Cells.Item(652,3).Value = "사회복지학부"
$ws.Cells.Item(652,4).Value2 = 20232317
$ws.Cells.Item(652,5).Value = "김은별"
$ws.Cells.Item(652,6).Value = "74:26"
$ws.Cells.Item(652,7).Value2 = 0.2
$ws.Cells.Item(652,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(652,9).Value = "952만 명"
$ws.Cells.Item(652,10).Value2 = 0.059
$ws.Cells.Item(652,11).Value = "중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다"
$ws.Cells.Item(652,12).Value = "Black"
$ws.Cells.Item(652,14).Value = "모름/무응답"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A652:L652")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcN = $ws.Range("N637")
$dstN = $ws.Range("N652")
$srcN.Copy()
$dstN.PasteSpecial(-4122)

# Row 653
$ws.Cells.Item(653,1).Value2 = 45192.804232407405
$ws.Cells.Item(653,2).Value = "kangsamy2@gmail.com"
$ws.Cells.Item(653,3).Value = "사회복지학부"
$ws.Cells.Item(653,4).Value2 = 20232302
$ws.Cells.Item(653,5).Value = "강새미"
$ws.Cells.Item(653,6).Value = "76:24"
$ws.Cells.Item(653,7).Value2 = 0.2
$ws.Cells.Item(653,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(653,9).Value = "779만 명"
$ws.Cells.Item(653,10).Value2 = 0.151
$ws.Cells.Item(653,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(653,12).Value = "Red"
$ws.Cells.Item(653,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$src = $ws.Range("A638:L638")
$dst = $ws.Range("A653:L653")
$src.Copy()
$dst.PasteSpecial(-4122)
$srcM = $ws.Range("M638")
$dstM = $ws.Range("M653")
$srcM.Copy()
$dstM.PasteSpecial(-4122)

$excel.CutCopyMode = 0
